# Updated data with DX explanations
#
# The WHODAS_SR question text in A16 gains a leading space before
# "In the past 30 days, for how many days were you totally unable
# to carry out your usual activities or work because of any health
# condition?"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A16").Value = " In the past 30 days, for how many days were you totally unable to carry out your usual activities or work because of any health condition?"
